$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list values (Price column D, Volume(1h) column E).
# D-column values that look like plain numbers (e.g. "1.001", "45.44") would
# otherwise be auto-converted to numeric cells by Excel's type inference, so
# we briefly force a text ("@") number format before assigning, then restore
# the "Normal" style so the cell's style index stays the same as before.
$updates = @(
    @{Row=2; D='27.753.87'; E='  -1.76%  '},
    @{Row=3; D='1.759.11'; E='  -2.06%  '},
    @{Row=4; E='  -0.05%  '},
    @{Row=5; D='325.44'; E='  -3.84%  '},
    @{Row=6; D='1.001'; E='  +0.01%  '},
    @{Row=7; D='0.4452'; E='  -2.54%  '},
    @{Row=8; D='0.3751'; E='  +0.38%  '},
    @{Row=9; D='45.44'; E='  +0.55%  '},
    @{Row=10; D='0.07541'; E='  -0.70%  '},
    @{Row=11; E='  -1.39%  '},
    @{Row=12; E='  -0.11%  '},
    @{Row=13; D='21.77'; E='  -2.51%  '},
    @{Row=14; D='6.214'; E='  -1.11%  '},
    @{Row=15; D='7.377'; E='  -1.26%  '},
    @{Row=16; D='1.761.80'; E='  -1.98%  '},
    @{Row=17; D='0.00001074'; E='  -1.37%  '},
    @{Row=18; D='88.01'; E='  +8.48%  '},
    @{Row=19; D='0.06227'; E='  -7.68%  '},
    @{Row=20; D='1.001'; E='  +0.04%  '},
    @{Row=21; D='17.34'; E='  -0.36%  '},
    @{Row=22; D='6.190'; E='  -2.81%  '},
    @{Row=23; D='0.5331'; E='  -3.80%  '},
    @{Row=24; D='27.782.06'; E='  -1.65%  '},
    @{Row=25; E='  -1.15%  '},
    @{Row=26; D='2.314'; E='  -4.36%  '},
    @{Row=27; D='20.69'; E='  +0.51%  '},
    @{Row=28; D='153.41'; E='  +0.72%  '},
    @{Row=29; D='2.369'; E='  +0.95%  '},
    @{Row=30; D='1.958.27'; E='  -2.24%  '},
    @{Row=31; D='128.41'; E='  -3.62%  '},
    @{Row=32; D='1.224'; E='  -0.67%  '},
    @{Row=33; D='0.09352'; E='  -0.88%  '},
    @{Row=34; D='5.753'; E='  -0.79%  '},
    @{Row=35; D='3.648'; E='  -9.50%  '},
    @{Row=36; E='  +5.48%  '},
    @{Row=37; D='0.02336'; E='  -0.22%  '},
    @{Row=38; D='0.2178'; E='  -7.16%  '},
    @{Row=39; E='  -2.64%  '},
    @{Row=40; D='0.6500'; E='  -1.22%  '},
    @{Row=41; D='5.089'; E='  -2.68%  '},
    @{Row=42; D='1.202'; E='  -0.75%  '},
    @{Row=43; D='8.010'; E='  -3.82%  '},
    @{Row=44; D='1.418'; E='  -4.38%  '},
    @{Row=45; E='  +0.00%  '},
    @{Row=46; D='13.81'; E='  -2.41%  '},
    @{Row=47; D='0.6015'; E='  -1.24%  '},
    @{Row=48; D='3.755'; E='  -2.53%  '},
    @{Row=49; D='126.58'; E='  -2.95%  '},
    @{Row=51; D='0.06907'; E='  -3.20%  '}
)

foreach ($item in $updates) {
    $r = $item.Row
    if ($item.ContainsKey('D')) {
        $cell = $ws.Range("D$r")
        $cell.NumberFormat = "@"
        $cell.Value = $item.D
        $cell.Style = "Normal"
    }
    if ($item.ContainsKey('E')) {
        $ws.Range("E$r").Value = $item.E
    }
}
